$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel's automatic
# date-pattern detection clobbering it (which would both turn the string
# into a date serial AND allocate a brand-new cell style). We stage the
# text as a formula result in a scratch cell, copy it, and paste-special
# "values only" into the destination -- that preserves the destination's
# existing style while landing a plain literal string.
function Set-LiteralText {
    param($range, [string]$text)

    $scratch = $ws.Range("Z100")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163, 0)
    $scratch.ClearContents()
}

# --- Bad Drivers table updates (rows 3-5) ---
$ws.Range("C3").Value = 3969
$ws.Range("D3").Value = 71.90000000000001

$ws.Range("C4").Value = 6245
$ws.Range("D4").Value = 94

$ws.Range("C5").Value = 10214

# --- Good Drivers table (rows 13-18): reordered by driver vintage, plus
#     refreshed sample counts / percentages / newly-known vintage dates ---
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001
Set-LiteralText $ws.Range("E13") "2024-11-10"

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B14").Value = 77849
$ws.Range("D14").Value = 99.90000000000001
Set-LiteralText $ws.Range("E14") "2021-08-18"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B15").Value = 34244
$ws.Range("D15").Value = 100
Set-LiteralText $ws.Range("E15") "2021-04-27"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100
Set-LiteralText $ws.Range("E16") "2020-08-05"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100
Set-LiteralText $ws.Range("E17") "2020-01-06"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B18").Value = 56018
$ws.Range("D18").Value = 100
Set-LiteralText $ws.Range("E18") "2019-12-14"
